# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds the "K" series. Re-write the recalculated values for rows 2-15.
$newK = @{
    2  = 6
    3  = 3
    4  = 7
    5  = 3
    6  = 4
    7  = 3
    8  = 6
    9  = 5
    10 = 10
    11 = 2
    12 = 4
    13 = 5
    14 = 4
    15 = 5
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
